# edit.ps1 — applies the "Trade #20 closed" commit to live_trading_results.xlsx
#
# Summary of the change:
#   - Summary sheet: populated with overall + per-strategy aggregate rows.
#   - leadlag sheet: Trade #1 (row 2) and Trade #2 (row 3) flip from OPEN to
#     CLOSED with exit price / P&L / exit-reason / duration filled in; a new
#     Trade #20 (row 19) is appended as OPEN; a few column widths widen.
#   - All Trades sheet: populated with the header row + the two now-closed
#     trades (mirrors leadlag rows 2-3).
#   - Comparison sheet: populated with the header row + the leadlag summary
#     stats row.
#
# Helper functions below force string values to stay literal text (the
# engine otherwise auto-detects numbers/dates/percentages from plain
# strings like "0.0%" or "2026-02-16" and silently converts them), while
# leaving the cell style untouched (ClearFormats after the write removes
# the transient "@" text format so cells keep the workbook's default style).

$wb = $excel.ActiveWorkbook

function Set-Text($ws, $addr, [string]$text) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

function Set-Num($ws, $addr, $num) {
    $ws.Range($addr).Value = $num
}

function Set-ColWidth($ws, $colIndex, $width) {
    # Excel's ColumnWidth property is offset by ~5/6 of a character from the
    # raw OOXML <col width="..."> value; subtract that offset so the saved
    # file ends up with the exact target width.
    $ws.Columns.Item($colIndex).ColumnWidth = $width - (5.0 / 6.0)
}

# ── Summary sheet ────────────────────────────────────────────────────────
$wsSummary = $wb.Worksheets.Item("Summary")

Set-Text $wsSummary "A1" "Metric"
Set-Text $wsSummary "B1" "Strategy"
Set-Text $wsSummary "C1" "Total Trades"
Set-Text $wsSummary "D1" "Win Rate"
Set-Text $wsSummary "E1" "Total P&L %"
Set-Text $wsSummary "F1" "Avg Trade"

Set-Text $wsSummary "A2" "OVERALL"
Set-Text $wsSummary "B2" "ALL COMBINED"
Set-Num  $wsSummary "C2" 2
Set-Text $wsSummary "D2" "0.0%"
Set-Text $wsSummary "E2" "-0.4739%"
Set-Text $wsSummary "F2" "-0.2369%"

Set-Text $wsSummary "A3" "STRATEGY"
Set-Text $wsSummary "B3" "leadlag"
Set-Num  $wsSummary "C3" 17
Set-Text $wsSummary "D3" "0.0%"
Set-Text $wsSummary "E3" "-0.4739%"
Set-Text $wsSummary "F3" "-0.0279%"

Set-ColWidth $wsSummary 1 10
Set-ColWidth $wsSummary 2 14
Set-ColWidth $wsSummary 3 14
Set-ColWidth $wsSummary 4 10
Set-ColWidth $wsSummary 5 13
Set-ColWidth $wsSummary 6 11

# ── leadlag sheet ────────────────────────────────────────────────────────
$wsLeadlag = $wb.Worksheets.Item("leadlag")

# Trade #1 (row 2): OPEN -> CLOSED
Set-Num  $wsLeadlag "G2" 69657.02966099999
Set-Text $wsLeadlag "H2" "CLOSED"
Set-Num  $wsLeadlag "I2" -0.1638
Set-Num  $wsLeadlag "J2" -1.64
Set-Text $wsLeadlag "M2" "time_exit_5min"
Set-Num  $wsLeadlag "N2" 5

# Trade #2 (row 3): OPEN -> CLOSED
Set-Num  $wsLeadlag "G3" 69709.80743099999
Set-Text $wsLeadlag "H3" "CLOSED"
Set-Num  $wsLeadlag "I3" -0.3101
Set-Num  $wsLeadlag "J3" -3.1
Set-Text $wsLeadlag "M3" "time_exit_5min"
Set-Num  $wsLeadlag "N3" 5

# Trade #20 (new row 19): OPEN
Set-Num  $wsLeadlag "A19" 20
Set-Text $wsLeadlag "B19" "2026-02-16"
Set-Text $wsLeadlag "C19" "21:25:29"
Set-Text $wsLeadlag "D19" "leadlag"
Set-Text $wsLeadlag "E19" "DOWN"
Set-Num  $wsLeadlag "F19" 69143.23
Set-Text $wsLeadlag "H19" "OPEN"
Set-Num  $wsLeadlag "I19" 0
Set-Num  $wsLeadlag "J19" 0
Set-Num  $wsLeadlag "K19" 0.75
Set-Text $wsLeadlag "L19" "Binance leading with -0.155% move"
Set-Num  $wsLeadlag "N19" 0

Set-ColWidth $wsLeadlag 7 14
Set-ColWidth $wsLeadlag 9 9
Set-ColWidth $wsLeadlag 13 16

# ── All Trades sheet ─────────────────────────────────────────────────────
$wsAllTrades = $wb.Worksheets.Item("All Trades")

Set-Text $wsAllTrades "A1" "Trade #"
Set-Text $wsAllTrades "B1" "Date"
Set-Text $wsAllTrades "C1" "Time"
Set-Text $wsAllTrades "D1" "Strategy"
Set-Text $wsAllTrades "E1" "Side"
Set-Text $wsAllTrades "F1" "Entry Price"
Set-Text $wsAllTrades "G1" "Exit Price"
Set-Text $wsAllTrades "H1" "Status"
Set-Text $wsAllTrades "I1" "P&L %"
Set-Text $wsAllTrades "J1" "P&L $"
Set-Text $wsAllTrades "K1" "Confidence"
Set-Text $wsAllTrades "L1" "Entry Reason"
Set-Text $wsAllTrades "M1" "Exit Reason"
Set-Text $wsAllTrades "N1" "Duration (min)"

Set-Num  $wsAllTrades "A2" 1
Set-Text $wsAllTrades "B2" "2026-02-16"
Set-Text $wsAllTrades "C2" "21:20:19"
Set-Text $wsAllTrades "D2" "leadlag"
Set-Text $wsAllTrades "E2" "DOWN"
Set-Num  $wsAllTrades "F2" 69543.125
Set-Num  $wsAllTrades "G2" 69657.02966099999
Set-Text $wsAllTrades "H2" "CLOSED"
Set-Num  $wsAllTrades "I2" -0.1638
Set-Num  $wsAllTrades "J2" -1.64
Set-Num  $wsAllTrades "K2" 0.7024
Set-Text $wsAllTrades "L2" "Binance leading with -0.070% move"
Set-Text $wsAllTrades "M2" "time_exit_5min"
Set-Num  $wsAllTrades "N2" 5

Set-Num  $wsAllTrades "A3" 2
Set-Text $wsAllTrades "B3" "2026-02-16"
Set-Text $wsAllTrades "C3" "21:20:25"
Set-Text $wsAllTrades "D3" "leadlag"
Set-Text $wsAllTrades "E3" "DOWN"
Set-Num  $wsAllTrades "F3" 69494.32000000001
Set-Num  $wsAllTrades "G3" 69709.80743099999
Set-Text $wsAllTrades "H3" "CLOSED"
Set-Num  $wsAllTrades "I3" -0.3101
Set-Num  $wsAllTrades "J3" -3.1
Set-Num  $wsAllTrades "K3" 0.75
Set-Text $wsAllTrades "L3" "Binance leading with -0.160% move"
Set-Text $wsAllTrades "M3" "time_exit_5min"
Set-Num  $wsAllTrades "N3" 5

Set-ColWidth $wsAllTrades 1 9
Set-ColWidth $wsAllTrades 2 12
Set-ColWidth $wsAllTrades 3 10
Set-ColWidth $wsAllTrades 4 10
Set-ColWidth $wsAllTrades 5 6
Set-ColWidth $wsAllTrades 6 13
Set-ColWidth $wsAllTrades 7 14
Set-ColWidth $wsAllTrades 8 8
Set-ColWidth $wsAllTrades 9 9
Set-ColWidth $wsAllTrades 10 7
Set-ColWidth $wsAllTrades 11 12
Set-ColWidth $wsAllTrades 12 35
Set-ColWidth $wsAllTrades 13 16
Set-ColWidth $wsAllTrades 14 16

# ── Comparison sheet ─────────────────────────────────────────────────────
$wsComparison = $wb.Worksheets.Item("Comparison")

Set-Text $wsComparison "A1" "Strategy"
Set-Text $wsComparison "B1" "Total Trades"
Set-Text $wsComparison "C1" "Win Rate"
Set-Text $wsComparison "D1" "Profit Factor"
Set-Text $wsComparison "E1" "Avg Win %"
Set-Text $wsComparison "F1" "Avg Loss %"
Set-Text $wsComparison "G1" "Win/Loss Ratio"
Set-Text $wsComparison "H1" "Max Drawdown"

Set-Text $wsComparison "A2" "leadlag"
Set-Num  $wsComparison "B2" 17
Set-Text $wsComparison "C2" "0.0%"
Set-Text $wsComparison "D2" "0.00"
Set-Text $wsComparison "E2" "+0.0000%"
Set-Text $wsComparison "F2" "-0.2369%"
Set-Text $wsComparison "G2" "0.00"
Set-Text $wsComparison "H2" "-0.3101%"

Set-ColWidth $wsComparison 1 10
Set-ColWidth $wsComparison 2 14
Set-ColWidth $wsComparison 3 10
Set-ColWidth $wsComparison 4 15
Set-ColWidth $wsComparison 5 11
Set-ColWidth $wsComparison 6 12
Set-ColWidth $wsComparison 7 16
Set-ColWidth $wsComparison 8 14

Write-Output "Edit complete."
